# CST-247 Burndown Chart — weekly actuals update
#
# Fills in this week's actual-hours entries for three team members
# (rows 5-7 of the Sprint burndown table) which ripple automatically
# through the "Actual Hours" totals (row 20) and "Ideal Burndown"
# tracker (row 21) via their existing SUM / IF formulas.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F5").Value = 5
$ws.Range("G6").Value = 6
$ws.Range("H7").Value = 5
$ws.Range("I7").Value = 3

# Recalculate so every dependent formula (row 20 totals, row 21 ideal
# burndown, and the downstream chart series) picks up the new inputs.
$excel.CalculateFullRebuild()

# Leave the cursor where the author left it when they saved.
$ws.Range("P12").Select() | Out-Null
